# Adds a new "to do" row (row 40) to the PLLF tracker table, motivated by a
# failure encountered with a mixed model. This mirrors what the author did in
# Excel: a new record was appended below the existing table, the table/
# autofilter range was grown to include it, and the worksheet selection was
# left on the empty row just below the new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles, wrap text, date number format, etc.) of the
# last existing row down into the new row before filling in values.
$ws.Range("A39:E39").Copy()
$ws.Range("A40:E40").PasteSpecial(-4122)

# Grow the table (ListObject) so its range/autoFilter covers the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:E40"))

# Populate the new row: problem, Raised by, Noted date. Action/date are left
# blank, same as in the source edit.
$ws.Range("A40").Value = "check whether old code worked for mixed; document that new code doesn't because neither offset nor constraint is allowed (despite the help file saying constraint is allowed), and that meglm can be tried"
$ws.Range("B40").Value = "Ian"
$ws.Range("C40").Value = 45981

# The long wrapped problem text needs a taller row.
$ws.Rows.Item(40).RowHeight = 58

# Leave the cursor on the next empty row, as in the authored workbook.
$ws.Range("A41").Select() | Out-Null
